$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A (the GENE column that was duplicated at the front of the
# table); this shifts columns B:F left to become the new A:E, matching the
# diff where the sheet's dimension changes from A1:F8 to A1:E8.
$ws.Columns("A").Delete()
